$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 "Modelo", matching the style of the existing header row
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Update numeric values in row 2
$ws.Range("B2").Value = 0.4202169370115176
$ws.Range("C2").Value = 0.9917138973304156
$ws.Range("D2").Value = 0.5082019190462719

# Add new cell F2 with the model description
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor())])"
